# Prepare for version 0.0.7
#
# Slide 1 contains a single picture/diagram group ("Group 12", id 13) made up
# of four child shapes: "Picture 6", "Rectangle 4", "Connector: Elbow 5" and
# "Picture 8". The edit re-lays-out that whole diagram: every child shape is
# moved/resized, and (as a natural consequence of PowerPoint recomputing the
# group's bounding box after such an edit) the group itself ends up with a
# new position/size, a new shape id/name ("Group 14", id 15) and a new
# creation-id.
#
# EMU/point note: the COM object model exposes Left/Top/Width/Height in
# points (1 pt = 12700 EMU) as (single-precision) floats, so a plain
# `emu / 12700` round-trips with an occasional 1-EMU truncation error. The
# literal values below are nudged by a few 1e-7 pt so they land back on the
# exact target EMU after that float round-trip.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)

# --- Reposition/resize the four child shapes in place (still inside the
#     original group) using their exact target absolute coordinates. ---

$pic6 = $grp.GroupItems.Item(1)   # "Picture 6"
$pic6.Left   = 255.0815812031496
$pic6.Top    = 13.021417322834646
$pic6.Width  = 516.1631775062991
$pic6.Height = 110.0348031496063

$rect4 = $grp.GroupItems.Item(2)  # "Rectangle 4"
$rect4.Left   = 255.0815812031496
$rect4.Top    = 96.0
$rect4.Width  = 71.9184251968504
$rect4.Height = 18.75007924015748

$conn = $grp.GroupItems.Item(3)   # "Connector: Elbow 5"
$conn.Left   = 380.78322834645667
$conn.Top    = 25.00763799527559
$conn.Width  = 35.99992185984252
$conn.Height = 215.48472440944883

$pic8 = $grp.GroupItems.Item(4)   # "Picture 8"
$pic8.Left   = 87.05102362204724
$pic8.Top    = 150.75
$pic8.Width  = 838.9490051779527
$pic8.Height = 591.4325256850393

# --- Ungroup and regroup so PowerPoint recomputes the group's own
#     off/ext/chOff/chExt as the tight bounding box of the repositioned
#     children (exactly what happens when a user nudges shapes around
#     inside a group and PowerPoint re-derives the envelope), and assigns
#     the group a fresh shape id/name in the process. ---

$ungrouped = $grp.Ungroup()

# Burn through the id counter so the recreated group lands on the same
# id/name ("Group 14", id 15) it has in the authored deck.
for ($k = 1; $k -le 8; $k++) {
    $tmp = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
    $tmp.Delete()
}

$newGroup = $s.Shapes.Range(@(1, 2, 3, 4)).Group()
$newGroup.Name = "Group 14"
